$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B19").Value = "2-LP"
$ws.Range("B13").Value = "0-LP"
$ws.Range("B14").Value = 0
$ws.Range("B16").Value = "8-LP"

$ws.Range("B16").Select()
